# edit.ps1 — reproduce the target commit:
#   * "deletePerson(p)" -> "deleteTask(t)" on the sequence-diagram slide
#   * cached "datetimeFigureOut" footer field text 2/6/2017 -> 3/16/17
#     on the slide master, every slide layout, and the notes master
#     (this is the auto-updating date placeholder getting re-cached).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Slide text: deletePerson(p) -> deleteTask(t)
#    The text lives in two runs inside one textbox so we touch each
#    run's characters individually to keep both runs (and their
#    purple-font formatting) intact, matching the original structure.
# ---------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $full = $shp.TextFrame.TextRange.Text
        if ($full -eq "deletePerson(p)") {
            $tr = $shp.TextFrame.TextRange
            # First run: "deletePerson" (12 chars) -> "deleteTask"
            $run1 = $tr.Characters(1, 12)
            $run1.Text = "deleteTask"
            # Second run now starts right after "deleteTask" (10 chars): "(p)" -> "(t)"
            $tr2 = $shp.TextFrame.TextRange
            $run2 = $tr2.Characters(11, 3)
            $run2.Text = "(t)"
        }
    }
}

# ---------------------------------------------------------------
# 2) Date placeholder re-cache: 2/6/2017 -> 3/16/17
#    Touch the "Date Placeholder" shape's text on the slide master,
#    every custom layout, and the notes master.
# ---------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "2/6/2017") {
                $shp.TextFrame.TextRange.Text = "3/16/17"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
